$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the file/case identifier in A2 with the new name
$ws.Range("A2").Value = "1007000/00012345/Pan Pes"

# Move the active selection to D4 (as recorded in the saved view state)
$ws.Range("D4").Select()
